# PacketFormats.xlsx update
# - Adds a new "My I'm alive packet timeout has been changed" row to the
#   "Controller to Server Packet Format" table (Data index 4), and renumbers
#   the Data index column in that table to be 0-based instead of 1-based.
# - Renumbers the Data index column to be 0-based in the "Server to
#   Controller Packet Format", "GWT to Server Packet Format (SSL)" and
#   "Server to GWT Packet Format (SSL)" tables as well.
# - Changes the "Separate for Data : |" entry to "Separate for Data :
#   Semicolon" (the delimiter character used for packet data changed from
#   pipe to semicolon).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 7 for the new "I'm alive timeout changed" packet
#    entry. This shifts every row from the old row 8 onward down by one,
#    which is exactly the shift seen for the rest of the sheet in the diff.
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = "My I'm alive packet timeout has been changed"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = "Controller Name"
$ws.Range("D7").Value = "Value (in second)"

# 2. "Controller to Server Packet Format" table (rows 3-6): Data index column
#    becomes 0-based (was 1,2,3,4 -> now 0,1,2,3).
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 3

# 3. "Server to Controller Packet Format" table (now rows 11-14 after the
#    insert above): Data index column becomes 0-based.
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 3

# 4. "GWT to Server Packet Format (SSL)" table (now row 18): Data index
#    column becomes 0-based.
$ws.Range("B18").Value = 0

# 5. "Server to GWT Packet Format (SSL)" table (now row 22): Data index
#    column becomes 0-based.
$ws.Range("B22").Value = 0

# 6. The trailing "Separate for Data : |" note (now row 26) becomes
#    "Separate for Data : Semicolon".
$ws.Range("A26").Value = "Separate for Data : Semicolon"

# 7. Match the author's final selection.
$ws.Range("C19").Select()
